$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest snapshot of coinranking.com prices/1h-volume figures for the
# cryptos list (GitHub Actions scheduled refresh). Every value below is
# stored as literal text in the sheet (Price/Volume columns are plain
# strings, not numbers/percentages) -- several "Price" figures look like
# plain numbers (e.g. "234.03"), so a leading apostrophe is used to stop
# Excel from auto-converting them to numeric values, exactly as typing
# them into Excel by hand would.
$updates = @(
    @{ Cell = 'D2'; Text = '43.486.60' }
    @{ Cell = 'E2'; Text = '  -1.05%  ' }
    @{ Cell = 'D3'; Text = '2.371.75' }
    @{ Cell = 'E3'; Text = '  +5.60%  ' }
    @{ Cell = 'E4'; Text = '  -0.10%  ' }
    @{ Cell = 'D5'; Text = '234.03' }
    @{ Cell = 'E5'; Text = '  +1.25%  ' }
    @{ Cell = 'D6'; Text = '0.644' }
    @{ Cell = 'E6'; Text = '  +0.18%  ' }
    @{ Cell = 'D7'; Text = '69.33' }
    @{ Cell = 'E7'; Text = '  +9.88%  ' }
    @{ Cell = 'E8'; Text = '  -0.01%  ' }
    @{ Cell = 'D9'; Text = '0.458' }
    @{ Cell = 'E10'; Text = '  -1.15%  ' }
    @{ Cell = 'D11'; Text = '57.36' }
    @{ Cell = 'E11'; Text = '  -0.18%  ' }
    @{ Cell = 'D12'; Text = '26.44' }
    @{ Cell = 'E12'; Text = '  +1.38%  ' }
    @{ Cell = 'D13'; Text = '2.725.69' }
    @{ Cell = 'E13'; Text = '  +5.78%  ' }
    @{ Cell = 'D14'; Text = '0.106' }
    @{ Cell = 'E14'; Text = '  +0.83%  ' }
    @{ Cell = 'D15'; Text = '15.67' }
    @{ Cell = 'E15'; Text = '  +1.27%  ' }
    @{ Cell = 'D16'; Text = '6.22' }
    @{ Cell = 'E16'; Text = '  +1.86%  ' }
    @{ Cell = 'D17'; Text = '0.852' }
    @{ Cell = 'E17'; Text = '  +2.91%  ' }
    @{ Cell = 'D18'; Text = '2.374.49' }
    @{ Cell = 'E18'; Text = '  +6.15%  ' }
    @{ Cell = 'D19'; Text = '43.474.01' }
    @{ Cell = 'E19'; Text = '  -0.76%  ' }
    @{ Cell = 'D20'; Text = '0.0₃0986' }
    @{ Cell = 'E20'; Text = '  -0.01%  ' }
    @{ Cell = 'D21'; Text = '6.33' }
    @{ Cell = 'E21'; Text = '  +4.25%  ' }
    @{ Cell = 'D22'; Text = '73.99' }
    @{ Cell = 'E22'; Text = '  +1.86%  ' }
    @{ Cell = 'D23'; Text = '247.85' }
    @{ Cell = 'E23'; Text = '  +0.12%  ' }
    @{ Cell = 'D24'; Text = '3.96' }
    @{ Cell = 'E24'; Text = '  +17.99%  ' }
    @{ Cell = 'E25'; Text = '  +0.03%  ' }
    @{ Cell = 'E26'; Text = '  +1.73%  ' }
    @{ Cell = 'E27'; Text = '  +2.40%  ' }
    @{ Cell = 'D28'; Text = '22.84' }
    @{ Cell = 'E28'; Text = '  +8.87%  ' }
    @{ Cell = 'D29'; Text = '9.96' }
    @{ Cell = 'E29'; Text = '  +1.58%  ' }
    @{ Cell = 'D30'; Text = '172.59' }
    @{ Cell = 'E30'; Text = '  +0.40%  ' }
    @{ Cell = 'D31'; Text = '1.55' }
    @{ Cell = 'E31'; Text = '  +9.86%  ' }
    @{ Cell = 'D32'; Text = '0.127' }
    @{ Cell = 'E32'; Text = '  -8.36%  ' }
    @{ Cell = 'E33'; Text = '  +1.41%  ' }
    @{ Cell = 'D34'; Text = '4.98' }
    @{ Cell = 'E34'; Text = '  +4.32%  ' }
    @{ Cell = 'E35'; Text = '  +0.62%  ' }
    @{ Cell = 'D36'; Text = '5.10' }
    @{ Cell = 'E36'; Text = '  +3.12%  ' }
    @{ Cell = 'D37'; Text = '6.59' }
    @{ Cell = 'E37'; Text = '  +2.88%  ' }
    @{ Cell = 'E38'; Text = '  +7.45%  ' }
    @{ Cell = 'D39'; Text = '3.62' }
    @{ Cell = 'E39'; Text = '  -0.98%  ' }
    @{ Cell = 'E40'; Text = '  +0.73%  ' }
    @{ Cell = 'D41'; Text = '8.98' }
    @{ Cell = 'E41'; Text = '  +6.69%  ' }
    @{ Cell = 'E42'; Text = '  +0.18%  ' }
    @{ Cell = 'D43'; Text = '18.45' }
    @{ Cell = 'E43'; Text = '  +8.02%  ' }
    @{ Cell = 'D44'; Text = '1.19' }
    @{ Cell = 'E44'; Text = '  +10.95%  ' }
    @{ Cell = 'E45'; Text = '  +2.15%  ' }
    @{ Cell = 'D46'; Text = '99.18' }
    @{ Cell = 'E46'; Text = '  +1.98%  ' }
    @{ Cell = 'D47'; Text = '4.49' }
    @{ Cell = 'E47'; Text = '  +4.31%  ' }
    @{ Cell = 'D48'; Text = '0.0952' }
    @{ Cell = 'E48'; Text = '  +1.13%  ' }
    @{ Cell = 'D49'; Text = '1.448.91' }
    @{ Cell = 'E49'; Text = '  +1.12%  ' }
    @{ Cell = 'D50'; Text = '2.598.65' }
    @{ Cell = 'E50'; Text = '  +6.04%  ' }
    @{ Cell = 'D51'; Text = '9.83' }
    @{ Cell = 'E51'; Text = '  +1.20%  ' }
)

foreach ($u in $updates) {
    $text = $u.Text
    $cell = $ws.Range($u.Cell)
    if ($text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.Value = "'" + $text
    } else {
        $cell.Value = $text
    }
}
